$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text while writing price values that look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.189.33"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.843.23"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "241.30"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").Value = "0.6877"
$ws.Range("E6").Value = "  -2.19%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "0.2998"
$ws.Range("E8").Value = "  -2.23%  "

$ws.Range("D9").Value = "0.07466"
$ws.Range("E9").Value = "  -3.35%  "

$ws.Range("D10").Value = "23.23"
$ws.Range("E10").Value = "  -1.67%  "

$ws.Range("D11").Value = "0.07664"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").Value = "1.842.91"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").Value = "5.062"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").Value = "0.6832"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "87.29"
$ws.Range("E15").Value = "  -6.52%  "

$ws.Range("D16").Value = "6.171"
$ws.Range("E16").Value = "  -6.33%  "

$ws.Range("D17").Value = "29.181.58"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "0.000008185"
$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("D19").Value = "2.084.70"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "229.15"
$ws.Range("E20").Value = "  -5.27%  "

$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "7.403"
$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "0.1449"
$ws.Range("E25").Value = "  -4.06%  "

$ws.Range("D26").Value = "159.45"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "8.770"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "18.10"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").Value = "1.517"
$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").Value = "1.198"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").Value = "0.05279"
$ws.Range("E33").Value = "  +3.23%  "

$ws.Range("D34").Value = "0.7607"
$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("D35").Value = "1.854"
$ws.Range("E35").Value = "  -2.21%  "

$ws.Range("D36").Value = "1.135"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").Value = "1.306.31"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("D41").Value = "0.9348"
$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("D43").Value = "105.17"
$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "1.986.65"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5192"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "64.87"
$ws.Range("E47").Value = "  +0.67%  "

$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("D49").Value = "9.498"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").Value = "1.772"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "0.07402"
$ws.Range("E51").Value = "  +16.96%  "

# Restore column D style to Normal (remove temporary text format) to match original styling
$ws.Range("D2:D51").Style = "Normal"
